# Reorders the per-period "Estado de Cuenta" detail rows (B16:J127) into
# ascending chronological order (1604 -> 2103), interleaving the two
# workers (LINETH PAOLA CASSAS OROZCO / IDEURANDO URIBE DELGADO) for each
# period, and refreshes the "Salario Basico" column (G) to the new value.
# "Valor Mora" (F) keeps travelling with its period.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @("1604","1609","1610","1611","1612","1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712","1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812","1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912","2001","2002","2003","2004","2005","2006","2007","2008","2009","2010","2011","2012","2101","2102","2103")

$newSalario = 781242

$docNums = @{ "LINETH" = "45557044"; "IDEURANDO" = "6212281" }
$docNames = @{ "LINETH" = "LINETH PAOLA CASSAS OROZCO"; "IDEURANDO" = "IDEURANDO URIBE DELGADO" }

$row = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $period = $periods[$i]

    if ($period -eq "1604") {
        $valorMora = 27600
    } elseif ($period -eq "2103") {
        $valorMora = 23958
    } elseif ([int]$period -le 1808) {
        $valorMora = 27578
    } else {
        $valorMora = 31249
    }

    foreach ($who in @("LINETH", "IDEURANDO")) {
        $ws.Cells.Item($row, 2).Value = "CC"
        $ws.Cells.Item($row, 3).Value = $docNums[$who]
        $ws.Cells.Item($row, 4).Value = $docNames[$who]
        $ws.Cells.Item($row, 5).Value = $period
        $ws.Cells.Item($row, 6).Value = $valorMora
        $ws.Cells.Item($row, 7).Value = $newSalario
        $row = $row + 1
    }
}
